# Update odds/statistics in row 2 of the sheet to reflect the latest values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.05
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 3.25

$ws.Range("AI2").Value = 19
$ws.Range("AK2").Value = 34

$ws.Range("AN2").Value = 4.5
$ws.Range("AO2").Value = 11
